$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column F for the 2022 data (shifts old "total" column F -> G)
$ws.Columns("F").Insert()

# New column header (2022) and updated "total" values for each publisher
$ws.Range("F1").Value = 2022

$ws.Range("F2").Value = 891555.52
$ws.Range("F3").Value = 484439.77
$ws.Range("F4").Value = 96717.22
$ws.Range("F5").Value = 251230.77
$ws.Range("F6").Value = 26631.98
$ws.Range("F7").Value = 9695.11
$ws.Range("F8").Value = 23695.11
$ws.Range("F9").Value = 9093.98
$ws.Range("F10").Value = 5218.1499999999996
$ws.Range("F11").Value = 12486.67

$ws.Range("G2").Value = 2521143.23
$ws.Range("G3").Value = 1695514.05
$ws.Range("G4").Value = 993709.6
$ws.Range("G5").Value = 791376.54
$ws.Range("G6").Value = 211459.15
$ws.Range("G7").Value = 83457.179999999993
$ws.Range("G8").Value = 87189.93
$ws.Range("G9").Value = 66272.539999999994
$ws.Range("G10").Value = 52907.05
$ws.Range("G11").Value = 46112.94

# Data cells (B2:G11) all share the same 2-decimal numeric format
$ws.Range("B2:G11").NumberFormat = "0.00"

# Rebuild the filler/background rows below the table: previously rows 14-23
# (B:F), now rows 13-29 (B:G) sharing the same number format as the data.
$ws.Range("B13:G29").NumberFormat = "0.00"

Write-Host "done"
